# Fine-tune the bus coordinates of the 29-bus GB transmission network
# (Geo_lon / Geo_lat columns on the "bus" sheet).

$wb = $excel.ActiveWorkbook
$busSheet = $wb.Worksheets.Item("bus")

$coords = @(
    @{Row=2;  H=-4.67;                 I=50.35},
    @{Row=3;  H=-1.75;                 I=50.1},
    @{Row=4;  H=-3.32;                 I=50.690000000000005},
    @{Row=5;  H=-0.44;                 I=50.52},
    @{Row=6;  H=1;                     I=50.95},
    @{Row=7;  H=-0.05;                 I=51.370000000000005},
    @{Row=8;  H=-2.0499999999999998;   I=51.2},
    @{Row=9;  H=-4.37;                 I=51.67},
    @{Row=10; H=-0.7;                  I=51.79},
    @{Row=11; H=-2.4;                  I=51.71},
    @{Row=12; H=0.35;                  I=52.22},
    @{Row=13; H=0.87;                  I=52.980000000000004},
    @{Row=14; H=-0.56999999999999995;  I=52.56},
    @{Row=15; H=-1.62;                 I=52.56},
    @{Row=16; H=-1.88;                 I=52.13},
    @{Row=17; H=-2.8;                  I=52.980000000000004},
    @{Row=18; H=-1.88;                 I=52.980000000000004},
    @{Row=19; H=-0.83;                 I=52.9},
    @{Row=20; H=-4.5;                  I=53.07},
    @{Row=21; H=-3.06;                 I=53.660000000000004},
    @{Row=22; H=-2.3588399999999998;   I=53.539050000000003},
    @{Row=23; H=-1.49;                 I=53.57},
    @{Row=24; H=-0.31;                 I=53.57},
    @{Row=25; H=-0.96;                 I=54},
    @{Row=26; H=-2.0099999999999998;   I=54.08},
    @{Row=27; H=-2.67;                 I=54.85},
    @{Row=28; H=-1.23;                 I=55.1},
    @{Row=29; H=-0.44;                 I=54.76},
    @{Row=30; H=0.21;                  I=54.08}
)

foreach ($entry in $coords) {
    $busSheet.Cells.Item($entry.Row, 8).Value = $entry.H
    $busSheet.Cells.Item($entry.Row, 9).Value = $entry.I
}

# View / window state tweaks that accompanied the edit.
$lineSheet = $wb.Worksheets.Item("line&trafo")

$busSheet.Select()
$busSheet.Range("S14").Select()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 12

$lineSheet.Select()
$lineSheet.Range("N23").Select()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1

$busSheet.Select()
$busSheet.Range("S14").Select()

$busSheet.PageSetup.PaperSize = 9
$busSheet.PageSetup.Orientation = 1
